$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Pre-format the new rows' Sucursal column (C16:C20) the same way as the
# existing table (right-aligned text, same style as C15) *before* writing
# any values into them. That way the "0xx" branch codes are stored as text
# straight away (no leading apostrophe needed) and no new style gets minted
# in styles.xml - the cells simply reuse the existing style index.
$ws.Range("C15").Copy()
$ws.Range("C16:C20").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Add the new User / Sucursal rows. The order the brand-new text values are
# first entered matters: it determines the order the new unique strings are
# appended to the shared string table, and we want it to land exactly the
# way it happened originally - row 19 first, then the "|" marker in D15,
# then rows 16-18.
$ws.Range("A19").Value = "F00263"
$ws.Range("C19").Value = "063"

$ws.Range("D15").Value = "|"

$ws.Range("A16").Value = "F00644"
$ws.Range("C16").Value = "026"

$ws.Range("A17").Value = "F01106"
$ws.Range("C17").Value = "006"

$ws.Range("A18").Value = "F02547"
$ws.Range("C18").Value = "089"

$ws.Activate()
$null = $ws.Range("J11").Select()
